$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) updates: force text type (values look numeric) while
# preserving the default (unstyled) cell format, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.739.35"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.029.04"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.26"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.51"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.022.56"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.57"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.51"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.533.26"
$ws.Range("D16").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.760.55"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.031.62"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.63"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.24"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.08"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.09"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("D27").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("D29").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.49"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("D35").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.91"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.36"
$ws.Range("D40").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.06"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.300"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.96"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.69"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0358"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.718.83"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.65"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.26"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.30"
$ws.Range("D51").ClearFormats()

# E-column (Volume/1h change) updates: plain strings with padding/%,
# always stored as text, no special handling required.
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  +11.18%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +6.64%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("E38").Value = "  +8.09%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +3.69%  "
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("E43").Value = "  +10.80%  "
$ws.Range("E44").Value = "  +5.21%  "
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("E51").Value = "  +3.25%  "
